$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "316"
$ws.Range("D9").Value = "1042987.10"

$ws.Range("C11").Value = "532"
$ws.Range("D11").Value = "3836269.21"

$ws.Range("C17").Value = "707"
$ws.Range("D17").Value = "6573657.56"

$ws.Range("C52").Value = "798"
$ws.Range("D52").Value = "5186204.95"

$ws.Range("C80").Value = "455"
$ws.Range("D80").Value = "1479752.96"

$ws.Range("C82").Value = "1283"
$ws.Range("D82").Value = "10292750.17"

$ws.Range("C94").Value = "270"
$ws.Range("D94").Value = "738250.00"

$ws.Range("C96").Value = "656"
$ws.Range("D96").Value = "4504105.68"

$ws.Range("C104").Value = "1705"
$ws.Range("D104").Value = "9801214.95"

$ws.Range("C106").Value = "1666"
$ws.Range("D106").Value = "9105315.02"
